# Edit script: rewrite the two list-item paragraphs with the new design-rationale text.
# Uses Range.InsertXML with a pkg:package payload so paragraph properties (pStyle/numPr)
# are preserved where the payload's <w:p> carries no <w:pPr>, and dropped where it doesn't
# (matching the target: item 1 keeps its ListParagraph numbering, the new paragraphs do not).

$d = $word.ActiveDocument

# --- Step 1: rewrite the first list item's runs in place (keeps its ListParagraph/numPr). ---
$p1 = $d.Paragraphs(3)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)
$xmlPara1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Have the driver run each tick of run rather than </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>world.run</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>) ticking forever.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xmlPara1)

# --- Step 2: drop the blank separator paragraph between item 1 and item 2. ---
$pBlank = $d.Paragraphs(4)
$pBlank.Range.Delete()

# --- Step 3: replace the second list item with three plain (non-list) paragraphs. ---
$p2 = $d.Paragraphs(4)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$xmlRest = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The reason why a new class called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GameStart</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> was created was because</w:t></w:r><w:r><w:t xml:space="preserve"> there was no way to intercept each tick of </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>World.run</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) that prevents a lot of new functionality to be added between turns and made it impossible to end the game early when a certain number of turns was reached. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The proposal is to have Application run each tick of run would mean that between each </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>turn</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>World.run</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) just ticks everything forever. Impossible to intercept each tick without editing world code. Rather it should run for each tick and the driver should call the function for each </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>tick</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xmlRest)
